$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.487.79'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '3.036.98'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.90'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.05'
$ws.Range('E6').Value = '  +3.69%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.033.62'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.68'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('E11').Value = '  -1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.488'
$ws.Range('E12').Value = '  +7.28%  '
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.86'
$ws.Range('E14').Value = '  +6.73%  '
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '66.372.64'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').Value = '3.540.20'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.24'
$ws.Range('E18').Value = '  +4.81%  '
$ws.Range('D19').Value = '3.036.40'
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.37'
$ws.Range('E20').Value = '  +18.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '469.43'
$ws.Range('E21').Value = '  +3.63%  '
$ws.Range('E22').Value = '  +3.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.39'
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.82'
$ws.Range('E25').Value = '  +4.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.27'
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.08'
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.23'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.44'
$ws.Range('E30').Value = '  +1.44%  '
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.119'
$ws.Range('E32').Value = '  +8.48%  '
$ws.Range('D33').Value = '0.0₃0992'
$ws.Range('E33').Value = '  -4.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.29'
$ws.Range('E34').Value = '  +3.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.87'
$ws.Range('E36').Value = '  +1.10%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.993'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.46'
$ws.Range('E38').Value = '  +10.41%  '
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.56'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('E43').Value = '  -3.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.65'
$ws.Range('E44').Value = '  +2.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0362'
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '380.51'
$ws.Range('E46').Value = '  -5.27%  '
$ws.Range('D47').Value = '2.712.19'
$ws.Range('E47').Value = '  -2.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.88'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.60'
$ws.Range('E50').Value = '  +3.17%  '
$ws.Range('E51').Value = '  +3.91%  '
